$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Financial Model")

# Adjust Q3 projected revenue growth multiplier from 1.15 to 1.25
$ws.Range("M3").Formula = "=I3*1.25"

# Align the YoY growth formula for M27 with the sequential pattern used by N27
$ws.Range("M27").Formula = "=M3/L3-1"

# Reflect the editor's navigation: they were last looking at the Main sheet
# around AH11, then moved to the Financial Model sheet near the cell they
# edited (M3), which becomes the active tab/selection on save.
$wsMain = $wb.Worksheets.Item("Main")
$wsMain.Activate()
$wsMain.Range("AH11").Select()

$ws.Activate()
$ws.Range("M3").Select()
